$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Check employee login credentials against Active Directory" (row 19) now
# scores 15/15 on the points-earned column (was blank/0).
$ws.Range("C19").Value = 15

# "Handle credit cards with REST calls (requests library)" (row 22) moves
# from 10/15 to a full 15/15.
$ws.Range("C22").Value = 15

# Scroll/select state: the sheet was left scrolled to row 3 with C23 as the
# active cell (previously topLeftCell A2 / selection B12).
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C23").Select()
